$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 1
$ws.Range("C4").Value = $false
$ws.Range("D4").Value = "haobaobei"
$ws.Range("E4").Value = "beibaobao"
$ws.Range("F4").Value = 1
$ws.Range("H4").Value = $true
$ws.Range("L4").Value = "Resource_1"

$ws.Range("I14").Select() | Out-Null
